$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: split the " Styles" run into " S" / "tyles" and drop
#    the "_GoBack" bookmark in between (Word stamps this where the cursor
#    last was when the document was saved).
# ---------------------------------------------------------------------------

# "_GoBack" already exists elsewhere in the document (next to "three") -
# remove it there first since Word only ever keeps a single "_GoBack".
$existing = $d.Bookmarks.Item("_GoBack")
$existing.Delete()

# Find the point right after "Naïve S" (i.e. immediately before "tyles")
# and drop a collapsed bookmark there.
$found = $d.Content
$found.Find.Execute("Naïve S", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $found.End
$insertionPoint = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $insertionPoint)

# ---------------------------------------------------------------------------
# 2) Heading 1 style: force a page break before every Heading 1 paragraph.
# ---------------------------------------------------------------------------
$heading1 = $d.Styles.Item("Heading 1")
$heading1.ParagraphFormat.PageBreakBefore = $true
